$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B,C,E,F,G,J,L,M,N,O across rows 2-25
# (row 25 has no N/O change - those stay as-is)
$newValues = @{
    2 = @{ "B"=1.135570874805808; "C"=0.263318201289195; "E"=0.2273233935621874; "F"=1.978202155727672; "G"=0.002456290033249817; "J"=0.05095900013149546; "L"=0.4284708600604716; "M"=0.3304236578472697; "N"=1.560381953623658; "O"=3.419704047225878 }
    3 = @{ "B"=1.056281248217857; "C"=0.2558495972362778; "E"=0.2283396101667172; "F"=1.978883395599311; "G"=0.002459119679845207; "J"=0.04990776632901728; "L"=0.4240151795748943; "M"=0.3165252378589969; "N"=1.575662556045689; "O"=3.427084142481505 }
    4 = @{ "B"=1.00788868674411; "C"=0.2512210974274893; "E"=0.2290203206618227; "F"=1.980309571153157; "G"=0.002460950956834709; "J"=0.04926013669120266; "L"=0.4214239446488506; "M"=0.3081030873420758; "N"=1.585607079256857; "O"=3.433641998563701 }
    5 = @{ "B"=0.9882428740015143; "C"=0.2493242874414676; "E"=0.2293120160341005; "F"=1.981144345943378; "G"=0.002461720889805921; "J"=0.04899569821027328; "L"=0.4204045096023563; "M"=0.3046992958999084; "N"=1.589800977599314; "O"=3.436823714098807 }
    6 = @{ "B"=0.9849852404970818; "C"=0.2490086824933826; "E"=0.2293613164217128; "F"=1.981298280835503; "G"=0.00246185016849863; "J"=0.04895175736150392; "L"=0.4202374432905955; "M"=0.3041358159392615; "N"=1.590505914798285; "O"=3.437382791485277 }
    7 = @{ "B"=1.007623432970718; "C"=0.2511955594197843; "E"=0.2290241966255948; "F"=1.980319802211298; "G"=0.002460961244189552; "J"=0.04925657247631676; "L"=0.4214100481324152; "M"=0.3080570676807284; "N"=1.585663066931676; "O"=3.433682846322171 }
    8 = @{ "B"=1.108172001233868; "C"=0.2607519506909313; "E"=0.2276620248134424; "F"=1.978227946751844; "G"=0.002457246256878546; "J"=0.05059699691782527; "L"=0.4269046299241097; "M"=0.3256084824789767; "N"=1.565534064036743; "O"=3.421827967843143 }
    9 = @{ "B"=1.307617582136572; "C"=0.2791496813469081; "E"=0.2254397725382038; "F"=1.982118577148995; "G"=0.002450702709923975; "J"=0.05320756294125317; "L"=0.4388204109861107; "M"=0.3609024156056577; "N"=1.5305187055774; "O"=3.414673888229203 }
    10 = @{ "B"=1.455491057759957; "C"=0.2924547811930154; "E"=0.2240790244077466; "F"=1.989847181378835; "G"=0.002446342710039818; "J"=0.05511367049409444; "L"=0.4482634560890091; "M"=0.387357052321093; "N"=1.507505842787531; "O"=3.419254237490406 }
    11 = @{ "B"=1.523045484020997; "C"=0.2984610790005604; "E"=0.2235186649118912; "F"=1.994420361431452; "G"=0.002444455453457222; "J"=0.05597804879336721; "L"=0.45270744359145; "M"=0.399503909406107; "N"=1.497624885061377; "O"=3.423479746935982 }
    12 = @{ "B"=1.548666688533842; "C"=0.3007287804891803; "E"=0.2233148756823429; "F"=1.996304071061644; "G"=0.002443754549109703; "J"=0.05630495604509633; "L"=0.4544114471008669; "M"=0.4041195668889372; "N"=1.493967695346079; "O"=3.425388219742587 }
    13 = @{ "B"=1.543146954223687; "C"=0.3002406923538103; "E"=0.2233583918800228; "F"=1.995891624378118; "G"=0.0024439048904255; "J"=0.05623456950019445; "L"=0.4540435208738529; "M"=0.4031247993614642; "N"=1.494751577691282; "O"=3.424963475615812 }
    14 = @{ "B"=1.525152565665849; "C"=0.2986477800868954; "E"=0.2235017307213809; "F"=1.994572290755556; "G"=0.002444397514144045; "J"=0.05600495207953315; "L"=0.4528472098044034; "M"=0.3998833250979743; "N"=1.497322312034697; "O"=3.42363057598638 }
    15 = @{ "B"=1.514135628234499; "C"=0.2976711928767486; "E"=0.2235906239000194; "F"=1.993783945431545; "G"=0.002444701050262808; "J"=0.05586425017584062; "L"=0.4521171857919484; "M"=0.397899891909006; "N"=1.498907967488599; "O"=3.422854305364524 }
    16 = @{ "B"=1.451081859861688; "C"=0.2920613160023038; "E"=0.2241168231043229; "F"=1.989569581735509; "G"=0.002446467975669971; "J"=0.05505712464991319; "L"=0.4479760012413436; "M"=0.386565467331998; "N"=1.508163412951454; "O"=3.419021221810937 }
    17 = @{ "B"=1.412472735136191; "C"=0.2886079181427874; "E"=0.224454631911442; "F"=1.987254955337704; "G"=0.002447576501864059; "J"=0.05456126637902159; "L"=0.445473385273857; "M"=0.3796407875013443; "N"=1.513991855873726; "O"=3.417218596306611 }
    18 = @{ "B"=1.390292790178421; "C"=0.2866172669576343; "E"=0.2246544527976724; "F"=1.986023177298677; "G"=0.002448223148840706; "J"=0.05427580720823855; "L"=0.4440479151052159; "M"=0.37566850106343; "N"=1.51739953959099; "O"=3.416383356473375 }
    19 = @{ "B"=1.38278771995806; "C"=0.2859425229729879; "E"=0.2247230579920121; "F"=1.98562321671389; "G"=0.002448443649072082; "J"=0.05417911260908426; "L"=0.443567679455839; "M"=0.3743253830729216; "N"=1.518562823573163; "O"=3.416135169563091 }
    20 = @{ "B"=1.416579957315378; "C"=0.2889759889850438; "E"=0.2244181002591645; "F"=1.987491050662754; "G"=0.002447457560940221; "J"=0.05461407785553618; "L"=0.4457383485289483; "M"=0.3803768357589732; "N"=1.513365683012104; "O"=3.417389623176376 }
    21 = @{ "B"=1.530436883890559; "C"=0.2991158411159915; "E"=0.223459400689018; "F"=1.994955688117827; "G"=0.00244425244585773; "J"=0.05607240770213906; "L"=0.4531980223371335; "M"=0.4008349951401016; "N"=1.496564931299389; "O"=3.424013709109488 }
    22 = @{ "B"=1.605080394689423; "C"=0.305703411499394; "E"=0.2228818235039025; "F"=2.000719865402331; "G"=0.002442237885040458; "J"=0.05702308650125332; "L"=0.4581966224921388; "M"=0.4142981730295716; "N"=1.486077288496368; "O"=3.430140482475537 }
    23 = @{ "B"=1.565221009056643; "C"=0.302191140481284; "E"=0.2231856138564847; "F"=1.997562418079397; "G"=0.002443305779879092; "J"=0.05651592072521083; "L"=0.4555175499997688; "M"=0.4071042397213489; "N"=1.49162966410362; "O"=3.426705912255841 }
    24 = @{ "B"=1.414723029252627; "C"=0.2888096004923852; "E"=0.2244345987425618; "F"=1.987384003792386; "G"=0.002447511305039396; "J"=0.05459020298013684; "L"=0.4456185171487874; "M"=0.3800440408968129; "N"=1.513648598724281; "O"=3.417311675491533 }
    25 = @{ "B"=1.253423702416683; "C"=0.274209569749047; "E"=0.2259930684711993; "F"=1.98021063106917; "G"=0.002452393995007549; "J"=0.05250335608004164; "L"=0.4354754593254597; "M"=0.4071042397213489 }
}

foreach ($row in $newValues.Keys) {
    $rowData = $newValues[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
